# Product Requisition - Shortcut.xlsx : "this added by last report 10-05-25"
#
# Updates a handful of quantity/rate inputs on Sheet1 (the dependent
# "Total Value" formulas - shared formula si="0" (=D*C) for most rows,
# and explicit D*C / C*D formulas for a couple of rows - recalculate
# automatically).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$ws.Activate()

# Row 9: Qty 20 -> 45, Rate 470 -> 300  (Total Value E9 recalculates 9400 -> 13500)
$ws.Range("C9").Value = 45
$ws.Range("D9").Value = 300

# Row 14: Qty 50 -> 25  (Total Value E14 recalculates 16150 -> 8075)
$ws.Range("C14").Value = 25

# Row 31: Qty 2000 -> 7000  (Total Value E31 recalculates 37260 -> 130410)
$ws.Range("C31").Value = 7000

# Row 32: Qty 500 -> blank  (Total Value E32 recalculates 13985 -> 0)
$ws.Range("C32").ClearContents()

# Row 43: Qty 231901 -> 189210  (Total Value E43 recalculates 223204.7125 -> 182114.625)
$ws.Range("C43").Value = 189210

# E46 (grand total, =SUM(E7:E45)) updates automatically through recalculation.

# Final selection: active cell E14, last touched/reviewed cell (alongside E9).
$ws.Range("E9").Select()
$ws.Range("E14").Select()
